$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (so "Registration No" becomes the new column B,
# and Name/Email/English/Maths/Physics/Computer Science shift one column to the right).
$ws.Columns.Item(2).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(1, 2).Value = "Registration No"

# Fill the new "Registration No" column for all data rows (2-16) with the constant value.
$ws.Range("B2:B16").Value = "THAUSCS027"

# Column widths: new column B (Registration No) gets width ~18.2 (closest representable
# value given the runtime's pixel-snapped column-width rounding); column D (Email, shifted
# from the old column C) keeps its original width of 22.23 automatically via the insert.
$ws.Columns.Item(2).ColumnWidth = 17.3

# Update the active selection to match the recorded view state.
$ws.Range("B16").Select()
